# Updated cryptos list data (prices + 1h volume %) pulled on
# Thu Sep  7 11:46:18 UTC 2023 by the scheduled GitHub Actions job.
# Column D (Price) and E (Volume(1h)) are plain text cells in the
# source sheet (right-aligned via leading/trailing spaces for E),
# so numeric-looking D values are forced back to text after the
# COM 'smart' Value assignment auto-detects them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.738.46'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '1.629.90'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.26'
$ws.Range("D5").ClearFormats()
$ws.Range("E6").Value = '  -0.79%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  -0.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0632'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.47'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.01%  '
$ws.Range("E11").Value = '  +0.83%  '
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("D13").Value = '1.855.22'
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("D14").Value = '1.630.40'
$ws.Range("E14").Value = '  -0.43%  '
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("E16").Value = '  -2.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.03'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").Value = '25.757.12'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.44'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.76'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.23%  '
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("E23").Value = '  +1.75%  '
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("E25").Value = '  +3.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.86'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.83%  '
$ws.Range("E27").Value = '  +2.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.85'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("E29").Value = '  -0.67%  '
$ws.Range("E30").Value = '  -0.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0489'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.33'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.14%  '
$ws.Range("E33").Value = '  -1.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.56'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.64%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.37'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.34%  '
$ws.Range("E36").Value = '  +0.37%  '
$ws.Range("D37").Value = '1.133.30'
$ws.Range("E37").Value = '  +1.85%  '
$ws.Range("E38").Value = '  -2.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.541'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.74%  '
$ws.Range("E40").Value = '  -1.40%  '
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.52'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.19'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.51'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.05%  '
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("D46").Value = '1.764.78'
$ws.Range("E46").Value = '  -0.24%  '
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.33'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.41%  '
$ws.Range("E49").Value = '  +0.93%  '
$ws.Range("E50").Value = '  +0.12%  '
$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.33'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -7.20%  '
